$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.333.55"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.274.69"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.967"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "2.618.94"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "2.271.51"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "42.308.85"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0847"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.77%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.224"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "1.705.52"
$ws.Range("E46").Value = "  +6.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("E49").Value = "  -5.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("E51").Value = "  -2.42%  "
